$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ids = @(251790, 251231, 252350, 243536, 243526, 243527, 243535, 243523, 251685, 252785, 252784, 252783, 251849, 252277, 252456)
$dates = @(45846.58333333334, 45846.58333333334, 45847.58333333334, 45846.58333333334, 45846.58333333334, 45846.58333333334, 45846.58333333334, 45846.58333333334, 45846.58333333334, 45846.58333333334, 45846.58333333334, 45846.58333333334, 45846.58333333334, 45847.58333333334, 45847.58333333334)

$numFmt = $ws.Range("B2").NumberFormat

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
    $ws.Cells.Item($row, 2).Value = $dates[$i]
    $ws.Cells.Item($row, 2).NumberFormat = $numFmt
}
